$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10 (Objetivos:): the value moves from the old "Propiciar ao aluno..." text
# to the responsavel docente string.
$ws.Range("B10").Value = "2143261 - André Luis Ferraz"
$ws.Range("C10").Value = "2143261 - André Luis Ferraz"

# The old row 13 (A empty, B/C held "2143261 - André Luis Ferraz") is removed entirely,
# shifting every following row up by one.
$ws.Rows.Item(13).Delete()

# After the shift, former row 14 ("Programa resumido:") is now row 13; its value
# changes to "Semestral".
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

# Former row 16 ("Programa:") is now row 15; its value changes to the activation date.
$ws.Range("B15").Value = "01/01/2018"
$ws.Range("C15").Value = "01/01/2018"

# Former row 19 ("Método:") is now row 18; its value changes to the responsavel docente string.
$ws.Range("B18").Value = "2143261 - André Luis Ferraz"
$ws.Range("C18").Value = "2143261 - André Luis Ferraz"

# Former row 20 ("Critério:") is now row 19; its value changes to the evaluation method text.
$ws.Range("B19").Value = "A avaliação será feita por meio de provas escritas."
$ws.Range("C19").Value = "A avaliação será feita por meio de provas escritas."

# Former row 21 ("Norma de recuperação:") is now row 20; its value changes to the
# final-grade criterion text.
$ws.Range("B20").Value = "A Nota final (NF) será calculada da seguinte maneira:NF = (P1 + 2*P2)/3Sendo que para P2 a matéria será cumulativa do semestre."
$ws.Range("C20").Value = "A Nota final (NF) será calculada da seguinte maneira:NF = (P1 + 2*P2)/3Sendo que para P2 a matéria será cumulativa do semestre."

# Former row 22 ("Bibliografia:") is now row 21; its value changes to the recovery-norm text.
$ws.Range("B21").Value = "A recuperação será feita por meio de uma prova escrita (PR) e a média de recuperação (MR) calculada pela fórmula: MR = (NF + PR)/2"
$ws.Range("C21").Value = "A recuperação será feita por meio de uma prova escrita (PR) e a média de recuperação (MR) calculada pela fórmula: MR = (NF + PR)/2"
